$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = "bauen"
$ws.Cells.Item(3, 2).Value = "kosten"
$ws.Cells.Item(4, 2).Value = "seufzen"
$ws.Cells.Item(5, 2).Value = "zögern"
$ws.Cells.Item(6, 2).Value = "fallen"
$ws.Cells.Item(7, 2).Value = "kichern"
$ws.Cells.Item(8, 2).Value = "fügen"
$ws.Cells.Item(9, 2).Value = "rasen"
$ws.Cells.Item(10, 2).Value = "bergen"
$ws.Cells.Item(11, 2).Value = "filmen"
$ws.Cells.Item(12, 2).Value = "nähen"
$ws.Cells.Item(13, 2).Value = "wandern"
$ws.Cells.Item(14, 2).Value = "plaudern"
$ws.Cells.Item(15, 2).Value = "ehren"
$ws.Cells.Item(16, 2).Value = "stürmen"
$ws.Cells.Item(17, 2).Value = "zünden"
$ws.Cells.Item(18, 2).Value = "schaden"
$ws.Cells.Item(19, 2).Value = "äußern"
$ws.Cells.Item(20, 2).Value = "lesen"
$ws.Cells.Item(21, 2).Value = "tropfen"
$ws.Cells.Item(22, 2).Value = "wagen"
$ws.Cells.Item(23, 2).Value = "sagen"
$ws.Cells.Item(24, 2).Value = "lieben"
$ws.Cells.Item(25, 2).Value = "quälen"
$ws.Cells.Item(26, 2).Value = "liegen"
$ws.Cells.Item(27, 2).Value = "ächzen"
$ws.Cells.Item(28, 2).Value = "wohnen"
$ws.Cells.Item(29, 2).Value = "brauchen"
$ws.Cells.Item(30, 2).Value = "sinken"
$ws.Cells.Item(31, 2).Value = "warten"
$ws.Cells.Item(32, 2).Value = "erben"
$ws.Cells.Item(33, 2).Value = "kümmern"
$ws.Cells.Item(34, 2).Value = "leeren"
$ws.Cells.Item(35, 2).Value = "fahren"
$ws.Cells.Item(36, 2).Value = "schwimmen"
$ws.Cells.Item(37, 2).Value = "wechseln"
$ws.Cells.Item(38, 2).Value = "werden"
$ws.Cells.Item(39, 2).Value = "zielen"
$ws.Cells.Item(40, 2).Value = "gelten"
$ws.Cells.Item(41, 2).Value = "fangen"
$ws.Cells.Item(42, 2).Value = "tauchen"
$ws.Cells.Item(43, 2).Value = "krachen"
$ws.Cells.Item(44, 2).Value = "fischen"
$ws.Cells.Item(45, 2).Value = "münzen"
$ws.Cells.Item(46, 2).Value = "mauern"
$ws.Cells.Item(47, 2).Value = "stören"
$ws.Cells.Item(48, 2).Value = "herrschen"
$ws.Cells.Item(49, 2).Value = "rauben"
$ws.Cells.Item(50, 2).Value = "bellen"
$ws.Cells.Item(51, 2).Value = "regeln"
$ws.Cells.Item(52, 2).Value = "ruhen"
$ws.Cells.Item(53, 2).Value = "spielen"
$ws.Cells.Item(54, 2).Value = "hupen"
$ws.Cells.Item(55, 2).Value = "wenden"
$ws.Cells.Item(56, 2).Value = "zeichnen"
$ws.Cells.Item(57, 2).Value = "handeln"
$ws.Cells.Item(58, 2).Value = "schlucken"
$ws.Cells.Item(59, 2).Value = "sterben"
$ws.Cells.Item(60, 2).Value = "schrecken"
$ws.Cells.Item(61, 2).Value = "sichern"
$ws.Cells.Item(62, 2).Value = "lockern"
$ws.Cells.Item(63, 2).Value = "trauen"
$ws.Cells.Item(64, 2).Value = "fehlen"
$ws.Cells.Item(65, 2).Value = "schwören"
$ws.Cells.Item(66, 2).Value = "betteln"
$ws.Cells.Item(67, 2).Value = "schwächen"
$ws.Cells.Item(68, 2).Value = "schenken"
$ws.Cells.Item(69, 2).Value = "schultern"
$ws.Cells.Item(70, 2).Value = "sperren"
$ws.Cells.Item(71, 2).Value = "wüten"
$ws.Cells.Item(72, 2).Value = "schulden"
$ws.Cells.Item(73, 2).Value = "streifen"
$ws.Cells.Item(74, 2).Value = "flehen"
$ws.Cells.Item(75, 2).Value = "tollen"
$ws.Cells.Item(76, 2).Value = "heben"
$ws.Cells.Item(77, 2).Value = "wirken"
$ws.Cells.Item(78, 2).Value = "ärgern"
$ws.Cells.Item(79, 2).Value = "holen"
$ws.Cells.Item(80, 2).Value = "knurren"
$ws.Cells.Item(81, 2).Value = "zeigen"
$ws.Cells.Item(82, 2).Value = "jagen"
$ws.Cells.Item(83, 2).Value = "saufen"
$ws.Cells.Item(84, 2).Value = "meinen"
$ws.Cells.Item(85, 2).Value = "schreiten"
$ws.Cells.Item(86, 2).Value = "altern"
$ws.Cells.Item(87, 2).Value = "geben"
$ws.Cells.Item(88, 2).Value = "schämen"
$ws.Cells.Item(89, 2).Value = "hoffen"
$ws.Cells.Item(90, 2).Value = "stecken"
$ws.Cells.Item(91, 2).Value = "kehren"
$ws.Cells.Item(92, 2).Value = "wellen"
$ws.Cells.Item(93, 2).Value = "stellen"
$ws.Cells.Item(94, 2).Value = "drücken"
$ws.Cells.Item(95, 2).Value = "platzen"
$ws.Cells.Item(96, 2).Value = "treiben"
$ws.Cells.Item(97, 2).Value = "wehtun"
$ws.Cells.Item(98, 2).Value = "dienen"
$ws.Cells.Item(99, 2).Value = "freuen"
$ws.Cells.Item(100, 2).Value = "helfen"
$ws.Cells.Item(101, 2).Value = "dauern"
$ws.Cells.Item(102, 2).Value = "scheitern"
$ws.Cells.Item(103, 2).Value = "führen"
$ws.Cells.Item(104, 2).Value = "folgen"
$ws.Cells.Item(105, 2).Value = "pfeifen"
$ws.Cells.Item(106, 2).Value = "graben"
$ws.Cells.Item(107, 2).Value = "fällen"
$ws.Cells.Item(108, 2).Value = "stürzen"
$ws.Cells.Item(109, 2).Value = "biegen"
$ws.Cells.Item(110, 2).Value = "landen"
$ws.Cells.Item(111, 2).Value = "hauen"
$ws.Cells.Item(112, 2).Value = "hören"
$ws.Cells.Item(113, 2).Value = "boxen"
$ws.Cells.Item(114, 2).Value = "heilen"
$ws.Cells.Item(115, 2).Value = "scheinen"
$ws.Cells.Item(116, 2).Value = "wahren"
$ws.Cells.Item(117, 2).Value = "kennen"
$ws.Cells.Item(118, 2).Value = "drohen"
$ws.Cells.Item(119, 2).Value = "feiern"
$ws.Cells.Item(120, 2).Value = "passen"
$ws.Cells.Item(121, 2).Value = "irren"
$ws.Cells.Item(122, 2).Value = "formen"
$ws.Cells.Item(123, 2).Value = "feuern"
$ws.Cells.Item(124, 2).Value = "bitten"
$ws.Cells.Item(125, 2).Value = "schwingen"
$ws.Cells.Item(126, 2).Value = "stammen"
$ws.Cells.Item(127, 2).Value = "duschen"
$ws.Cells.Item(128, 2).Value = "mögen"
$ws.Cells.Item(129, 2).Value = "bilden"
$ws.Cells.Item(130, 2).Value = "wundern"
$ws.Cells.Item(131, 2).Value = "hindern"
$ws.Cells.Item(132, 2).Value = "schütteln"
$ws.Cells.Item(133, 2).Value = "warnen"
$ws.Cells.Item(134, 2).Value = "spüren"
$ws.Cells.Item(135, 2).Value = "bleiben"
$ws.Cells.Item(136, 2).Value = "siegen"
$ws.Cells.Item(137, 2).Value = "heulen"
$ws.Cells.Item(138, 2).Value = "prügeln"
$ws.Cells.Item(139, 2).Value = "liefern"
$ws.Cells.Item(140, 2).Value = "klingen"
$ws.Cells.Item(141, 2).Value = "töten"
$ws.Cells.Item(142, 2).Value = "dringen"
$ws.Cells.Item(143, 2).Value = "pissen"
$ws.Cells.Item(144, 2).Value = "drehen"
$ws.Cells.Item(145, 2).Value = "schmecken"
$ws.Cells.Item(146, 2).Value = "stoßen"
$ws.Cells.Item(147, 2).Value = "malen"
$ws.Cells.Item(148, 2).Value = "suchen"
$ws.Cells.Item(149, 2).Value = "rufen"
$ws.Cells.Item(150, 2).Value = "kämpfen"
$ws.Cells.Item(151, 2).Value = "enden"
$ws.Cells.Item(152, 2).Value = "gründen"
$ws.Cells.Item(153, 2).Value = "grüßen"
$ws.Cells.Item(154, 2).Value = "treten"
$ws.Cells.Item(155, 2).Value = "fließen"
$ws.Cells.Item(156, 2).Value = "ändern"
$ws.Cells.Item(157, 2).Value = "streichen"
$ws.Cells.Item(158, 2).Value = "lächeln"
$ws.Cells.Item(159, 2).Value = "räumen"
$ws.Cells.Item(160, 2).Value = "garen"
$ws.Cells.Item(161, 2).Value = "spinnen"
$ws.Cells.Item(162, 2).Value = "knarren"
$ws.Cells.Item(163, 2).Value = "arten"
$ws.Cells.Item(164, 2).Value = "planen"
$ws.Cells.Item(165, 2).Value = "sorgen"
$ws.Cells.Item(166, 2).Value = "trennen"
$ws.Cells.Item(167, 2).Value = "schlafen"
$ws.Cells.Item(168, 2).Value = "wachsen"
$ws.Cells.Item(169, 2).Value = "retten"
$ws.Cells.Item(170, 2).Value = "stärken"
$ws.Cells.Item(171, 2).Value = "weichen"
$ws.Cells.Item(172, 2).Value = "läuten"
$ws.Cells.Item(173, 2).Value = "klettern"
$ws.Cells.Item(174, 2).Value = "reizen"
$ws.Cells.Item(175, 2).Value = "decken"
$ws.Cells.Item(176, 2).Value = "flüchten"
$ws.Cells.Item(177, 2).Value = "jubeln"
$ws.Cells.Item(178, 2).Value = "werfen"
$ws.Cells.Item(179, 2).Value = "gnaden"
$ws.Cells.Item(180, 2).Value = "locken"
$ws.Cells.Item(181, 2).Value = "lügen"
$ws.Cells.Item(182, 2).Value = "loben"
$ws.Cells.Item(183, 2).Value = "leihen"
$ws.Cells.Item(184, 2).Value = "bluten"
$ws.Cells.Item(185, 2).Value = "sprengen"
$ws.Cells.Item(186, 2).Value = "greifen"
$ws.Cells.Item(187, 2).Value = "machen"
$ws.Cells.Item(188, 2).Value = "runden"
$ws.Cells.Item(189, 2).Value = "pflanzen"
$ws.Cells.Item(190, 2).Value = "klagen"
$ws.Cells.Item(191, 2).Value = "stehlen"
$ws.Cells.Item(192, 2).Value = "achten"
$ws.Cells.Item(193, 2).Value = "backen"
